$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.711.32'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '3.017.33'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  -0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '587.06'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.42%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '148.61'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("E7").Value = '  +0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.08%  '
$ws.Range("D9").Value = '3.015.51'
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("E10").Value = '  -2.39%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.84'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.90%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.463'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.22%  '
$ws.Range("E13").Value = '  -1.45%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '34.85'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.25%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.124'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '3.520.54'
$ws.Range("E16").Value = '  -0.71%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '7.12'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '62.676.87'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Value = '3.018.27'
$ws.Range("E19").Value = '  -1.04%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '461.20'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.90%  '
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("E22").Value = '  -1.63%  '
$ws.Range("E23").Value = '  -0.44%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '81.80'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("E25").Value = '  -7.71%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '12.37'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.85%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.03'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -5.94%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("E29").Value = '  -0.56%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.10%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.06'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -4.41%  '
$ws.Range("E32").Value = '  -4.04%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '28.10'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("D35").Value = '0.0₃0822'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("E36").Value = '  -2.23%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.79'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.22%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.13'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.82%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '50.43'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.08%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '9.19'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("E41").Value = '  -10.33%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.123'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +8.16%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '395.51'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -9.13%  '
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("E45").Value = '  -6.03%  '
$ws.Range("D46").Value = '2.747.43'
$ws.Range("E46").Value = '  -2.79%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '37.44'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.24%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '129.45'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.22'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
